# Updates cryptos list values (prices / 1h volume %) and fixes the
# swapped Polkadot / WrappedEther rows (12 and 13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, regardless of whether it "looks"
# numeric (e.g. "323.97", "1.000"), without leaving the cell's style
# changed (original workbook uses the default/general style on these
# data cells).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.101.25"
Set-TextValue $ws.Range("E2") "  -2.00%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.837.67"
Set-TextValue $ws.Range("E3") "  -0.77%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.20%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "323.97"
Set-TextValue $ws.Range("E5") "  -3.44%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "1.000"
Set-TextValue $ws.Range("E6") "  -0.25%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.4643"
Set-TextValue $ws.Range("E7") "  -0.38%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3874"
Set-TextValue $ws.Range("E8") "  -0.95%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("E9") "  -0.54%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "0.9631"
Set-TextValue $ws.Range("E10") "  -2.31%  "

# Row 11 - Solana
Set-TextValue $ws.Range("D11") "21.97"
Set-TextValue $ws.Range("E11") "  -1.59%  "

# Row 12 - now WrappedEther (was Polkadot)
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.850.59"
Set-TextValue $ws.Range("E12") "  -1.68%  "

# Row 13 - now Polkadot (was WrappedEther)
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "5.695"
Set-TextValue $ws.Range("E13") "  -2.72%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "6.924"
Set-TextValue $ws.Range("E14") "  -1.31%  "

# Row 15 - TRON
Set-TextValue $ws.Range("D15") "0.06818"
Set-TextValue $ws.Range("E15") "  -0.48%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "88.11"

# Row 17 - BinanceUSD
Set-TextValue $ws.Range("E17") "  -0.27%  "

# Row 18 - ShibaInu
Set-TextValue $ws.Range("E18") "  -1.10%  "

# Row 19 - Avalanche
Set-TextValue $ws.Range("E19") "  -2.47%  "

# Row 20 - Dai
Set-TextValue $ws.Range("E20") "  -0.22%  "

# Row 21 - WrappedBTC
Set-TextValue $ws.Range("D21") "28.105.17"
Set-TextValue $ws.Range("E21") "  -2.01%  "

# Row 22 - Uniswap
Set-TextValue $ws.Range("D22") "5.329"
Set-TextValue $ws.Range("E22") "  -1.44%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("E23") "  -2.39%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "2.099"
Set-TextValue $ws.Range("E24") "  -1.95%  "

# Row 25 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D25") "2.015.23"
Set-TextValue $ws.Range("E25") "  -4.56%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "154.80"

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "19.23"
Set-TextValue $ws.Range("E27") "  -1.33%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "5.687"
Set-TextValue $ws.Range("E28") "  -6.03%  "

# Row 29 - LidoDAOToken
Set-TextValue $ws.Range("D29") "1.966"
Set-TextValue $ws.Range("E29") "  -3.11%  "

# Row 30 - BitcoinCash
Set-TextValue $ws.Range("D30") "118.28"
Set-TextValue $ws.Range("E30") "  +0.38%  "

# Row 31 - ImmutableX
Set-TextValue $ws.Range("D31") "0.9380"
Set-TextValue $ws.Range("E31") "  -4.06%  "

# Row 32 - Stellar
Set-TextValue $ws.Range("D32") "0.09253"
Set-TextValue $ws.Range("E32") "  -1.96%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.278"
Set-TextValue $ws.Range("E33") "  -1.89%  "

# Row 34 - ARBITRUM
Set-TextValue $ws.Range("D34") "1.323"
Set-TextValue $ws.Range("E34") "  -2.19%  "

# Row 35 - HuobiToken
Set-TextValue $ws.Range("D35") "3.304"
Set-TextValue $ws.Range("E35") "  -5.20%  "

# Row 36 - Hedera
Set-TextValue $ws.Range("E36") "  -4.75%  "

# Row 37 - VeChain
Set-TextValue $ws.Range("D37") "0.02133"
Set-TextValue $ws.Range("E37") "  -3.08%  "

# Row 38 - TrustWalletToken
Set-TextValue $ws.Range("D38") "1.147"
Set-TextValue $ws.Range("E38") "  -1.41%  "

# Row 39 - FraxShare
Set-TextValue $ws.Range("D39") "7.782"
Set-TextValue $ws.Range("E39") "  +1.73%  "

# Row 40 - TheSandbox
Set-TextValue $ws.Range("D40") "0.5607"
Set-TextValue $ws.Range("E40") "  -2.22%  "

# Row 41 - Aptos
Set-TextValue $ws.Range("D41") "9.912"
Set-TextValue $ws.Range("E41") "  -3.11%  "

# Row 42 - Algorand
Set-TextValue $ws.Range("D42") "0.1765"
Set-TextValue $ws.Range("E42") "  -2.14%  "

# Row 43 - Cronos
Set-TextValue $ws.Range("D43") "0.07260"
Set-TextValue $ws.Range("E43") "  +1.60%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "11.70"
Set-TextValue $ws.Range("E44") "  -0.31%  "

# Row 45 - Decentraland
Set-TextValue $ws.Range("D45") "0.5278"
Set-TextValue $ws.Range("E45") "  -2.29%  "

# Row 46 - RenderToken
Set-TextValue $ws.Range("D46") "2.138"
Set-TextValue $ws.Range("E46") "  -10.36%  "

# Row 47 - WEMIXToken
Set-TextValue $ws.Range("D47") "1.136"
Set-TextValue $ws.Range("E47") "  -8.95%  "

# Row 48 - NEARProtocol
Set-TextValue $ws.Range("D48") "1.831"
Set-TextValue $ws.Range("E48") "  -4.12%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "112.89"
Set-TextValue $ws.Range("E49") "  -1.79%  "

# Row 50 - EOS
Set-TextValue $ws.Range("E50") "  +0.66%  "

# Row 51 - PaxDollar
Set-TextValue $ws.Range("E51") "  -0.19%  "
